$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, F, G, J, L across rows 2-25
# corresponding to the updated simulation results for the 380 kV case.
$newValues = @{
    2 = @(3.57734016330113, 0.1900693106435938, 0.03795911633739024, 6.362963440239838, 0.002662275072514337, 0.1988445409129653, 0.3795162915198773)
    3 = @(3.519276965316749, 0.1746043628738505, 0.03323333479210078, 6.195661185158173, 0.002669449000711936, 0.1973024380260995, 0.3802619368313387)
    4 = @(3.486581486025784, 0.1653137460002654, 0.03032715968228672, 6.094701191362333, 0.002674079829847691, 0.1963732021675249, 0.3809843144409797)
    5 = @(3.473998879121666, 0.1615786112528212, 0.02914147345811102, 6.0539951080126, 0.002676023983039921, 0.1959988655154881, 0.3813452209927704)
    6 = @(3.471954244787923, 0.1609614481247377, 0.028944499314008, 6.047261999480043, 0.002676350260563832, 0.1959369661855739, 0.3814091679742191)
    7 = @(3.486408794771762, 0.1652631675722489, 0.03031117509718229, 6.094150459928045, 0.002674105817879153, 0.1963681363092959, 0.3809889123500625)
    8 = @(3.55670520511535, 0.1846940566226181, 0.03633044936468366, 6.304906735861579, 0.002664701868473273, 0.1983091076847714, 0.3797184856477003)
    9 = @(3.718121456510914, 0.2244566125889094, 0.04811044471854586, 6.732569993043739, 0.002648043904094772, 0.2022597357078055, 0.3793267011677273)
    10 = @(3.851264216054801, 0.2547305862863709, 0.0567682092968198, 7.056086106862097, 0.002636878142887692, 0.2052569247861484, 0.380320642572002)
    11 = @(3.915036812405901, 0.2687444050382055, 0.0607113467977598, 7.205406551331407, 0.002632028456622893, 0.2066424449387796, 0.3810517236369435)
    12 = @(3.939650013702533, 0.2740867443575326, 0.06220548396885306, 7.262269135165582, 0.002630224797699286, 0.2071703960006488, 0.3813687173037152)
    13 = @(3.934328439335502, 0.2729345791649109, 0.06188364769623433, 7.250008467674206, 0.00263061179228119, 0.2070565444565204, 0.3812986607110247)
    14 = @(3.91705244044573, 0.2691832041338955, 0.06083424955237149, 7.210078242697421, 0.002631879411971607, 0.2066858133311484, 0.3810769980017312)
    15 = @(3.906530896865092, 0.266890040397584, 0.06019159594877976, 7.185661526320359, 0.002632660132804537, 0.2064591608324022, 0.3809464530630038)
    16 = @(3.847161278044155, 0.253819689062027, 0.05651063430813963, 7.046371705154286, 0.002637199684560922, 0.2051668314540649, 0.3802784807215289)
    17 = @(3.811562977928986, 0.2458640098596732, 0.05425388255082453, 6.961479146624782, 0.002640043224967079, 0.2043797625358508, 0.3799401682656196)
    18 = @(3.791389288161383, 0.2413108457389797, 0.05295628167309019, 6.912853151307019, 0.002641700384955934, 0.2039291333347819, 0.3797718315754111)
    19 = @(3.784610519140699, 0.2397731026617009, 0.05251700036967577, 6.896423645255226, 0.002642265192613837, 0.2037769106058249, 0.3797193433300237)
    20 = @(3.815321251268529, 0.2467085470322559, 0.0544940717963982, 6.970495139038377, 0.002639738288164311, 0.2044633321341962, 0.3799734649926592)
    21 = @(3.922114205492619, 0.2702841009359247, 0.06114245498466175, 7.221798011650208, 0.002631506192584809, 0.2067946160283327, 0.3811410157331494)
    22 = @(3.994615029568536, 0.2859000414834156, 0.06549326753861351, 7.387897710177469, 0.002626317197578915, 0.208337441969519, 0.382138120269417)
    23 = @(3.955671407517684, 0.2775462206184613, 0.06317053985550558, 7.299074150663614, 0.00262906924226432, 0.2075122137681689, 0.3815845173030255)
    24 = @(3.813621225831412, 0.2463266674806448, 0.05438548277155064, 6.966418449285783, 0.002639876080222522, 0.20442554451094, 0.3799583300513945)
    25 = @(3.671912626523863, 0.2135171956242345, 0.04492439081640498, 6.61528074887363, 0.002652360896223406, 0.2011748467198728, 0.3792077510134959)
}

$columns = @("B", "C", "D", "F", "G", "J", "L")

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "{0}{1}" -f $columns[$i], $row
        $ws.Range($cellRef).Value = $vals[$i]
    }
}
